# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The "Periodo Mora" (EC database) is rolled forward from period 2507/2506
# to a single consolidated period 2508, and the now-duplicate worker rows
# are removed, leaving one row per worker. The summary figures (Valor Mora
# total and Cant. Periodos) are updated to match the new, smaller data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (YAMELIS MARIA CABEZA LLERENA, CC 30854190) keeps its single
# period row, but the period itself rolls from 2507 -> 2508.
$ws.Range("E16").Value = "2508"

# Rows 17 & 18 were the extra period rows for this statement (YAMELIS'
# 2506 period, and CAROLINA's 2507 period) - both are now redundant and
# are removed entirely, shifting row 19 (CAROLINA, period 2506) up into
# row 17's position.
$ws.Rows("17:18").Delete()

# The surviving row (now row 17, CAROLINA PATRICIA LOPEZ BARBOSA, CC
# 1007978222) also rolls its period from 2506 -> 2508.
$ws.Range("E17").Value = "2508"

# "VALOR MORA" total now reflects only the two remaining rows.
$ws.Range("E11").Value = 136200

# "Cant. Periodos" drops from 2 to 1 since the data now covers a single
# consolidated period (2508).
$ws.Range("F13").Value = 1
